# Time_measurement.xlsx -- "tm final push +timestemp+date"
#
# The sheet's old sample data / formatting is replaced wholesale with a
# fresh, default-styled sheet containing the new timing measurements.
#
# We rebuild the worksheet from scratch (rather than editing cells in
# place) so none of the old custom fonts / column widths / row heights /
# page setup survive -- the target sheet uses Excel's plain defaults.

$wb = $excel.ActiveWorkbook

# Remember the sheet we need to replace, then add a brand-new blank sheet
# and get rid of the old one so no leftover styling/column sizing remains.
$oldSheet = $wb.ActiveSheet
$oldName = $oldSheet.Name

$ws = $wb.Worksheets.Add()
$newIndex = $ws.Index

# Locate the original sheet again (its Index shifts once the new sheet is
# inserted in front of it) and drop it, then reclaim its name.
foreach ($s in $wb.Worksheets) {
    if ($s.Index -ne $newIndex) {
        $null = $s.Delete()
        break
    }
}
$ws.Name = $oldName

$data = @(
    @("socorro", "8.344650268554688e-06"),
    @("socorro", "9.298324584960938e-06"),
    @("socorro", "8.58306884765625e-06"),
    @("socorro", "7.390975952148438e-06"),
    @("socorro", "7.867813110351562e-06"),
    @("socorro", "7.152557373046875e-07"),
    @("fundos fundos fundos", "4.76837158203125e-07"),
    @("fundos fundos fundos", "4.76837158203125e-07"),
    @("fundos fundos fundos", "7.152557373046875e-07"),
    @("socorro", "9.5367431640625e-07"),
    @("socorro", "9.5367431640625e-07"),
    @("socorro", "7.152557373046875e-07"),
    @("socorro", "4.76837158203125e-07")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = [double]$data[$i][1]
}
